# "Correction to size of ball formula."
#
# The D column ("size of ball") formula was shared only across D4:D41, and
# used MAX($C$4:$C$41) as the normalising max, while D42 held a hard-coded
# placeholder value of 20 instead of the real formula result. This extends
# the shared formula (and the MAX range it references) down to row 42 on
# all three Raw sheets, so D42 now computes its value the same way as every
# other row in the column.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Raw1", "Raw2", "Raw3")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Re-fill D4:D42 with the corrected formula. Excel adjusts the relative
    # C4 reference per-row automatically while the $C$4:$C$42 absolute range
    # stays fixed, exactly like the old D4:D41/$C$4:$C$41 shared formula did.
    $ws.Range("D4:D42").Formula = "=((MAX(`$C`$4:`$C`$42)/5+2)-C4/5)"
}

# Matches the saved selection state of the corrected workbook.
$wb.Worksheets.Item("Raw2").Range("D23").Select()
$wb.Worksheets.Item("Raw3").Range("D23").Select()

$wb.Worksheets.Item("Raw1").Activate()
